# Apply the price/volume refresh described by the commit diff.
# Numeric-looking "Price" strings (column D) are written with a leading
# apostrophe so Excel keeps them as literal text (matching the original
# inline-string formatting, e.g. trailing zeros like "0.4280") instead of
# silently converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.621.96"
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("D3").Value = "1.847.40"
$ws.Range("E3").Value = "  +0.07%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E5").Value = "  -0.33%  "
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").Value = "'0.4280"
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("D8").Value = "'0.3634"
$ws.Range("E8").Value = "  -0.34%  "
$ws.Range("D9").Value = "'44.94"
$ws.Range("E9").Value = "  +2.64%  "
$ws.Range("E10").Value = "  +1.19%  "
$ws.Range("D11").Value = "'0.8754"
$ws.Range("E11").Value = "  -2.52%  "
$ws.Range("D12").Value = "'20.68"
$ws.Range("E12").Value = "  +0.35%  "
$ws.Range("D13").Value = "1.858.51"
$ws.Range("E13").Value = "  +1.06%  "
$ws.Range("D14").Value = "'5.327"
$ws.Range("E14").Value = "  -0.82%  "
$ws.Range("D15").Value = "'6.518"
$ws.Range("E15").Value = "  -0.76%  "
$ws.Range("D16").Value = "'0.06910"
$ws.Range("E16").Value = "  +0.87%  "
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("D18").Value = "'79.91"
$ws.Range("E18").Value = "  +2.88%  "
$ws.Range("D19").Value = "'0.000009014"
$ws.Range("E19").Value = "  +1.26%  "
$ws.Range("D20").Value = "'1.001"
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("D21").Value = "'15.32"
$ws.Range("E21").Value = "  -0.60%  "
$ws.Range("D22").Value = "27.649.98"
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").Value = "'4.957"
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("D24").Value = "'10.39"
$ws.Range("E24").Value = "  -2.53%  "
$ws.Range("D25").Value = "2.095.14"
$ws.Range("E25").Value = "  +2.12%  "
$ws.Range("D26").Value = "'1.989"
$ws.Range("E26").Value = "  -2.89%  "
$ws.Range("D27").Value = "'155.04"
$ws.Range("E27").Value = "  +1.15%  "
$ws.Range("D28").Value = "'18.75"
$ws.Range("E28").Value = "  +2.63%  "
$ws.Range("D29").Value = "'121.37"
$ws.Range("E29").Value = "  +8.82%  "
$ws.Range("D30").Value = "'5.293"
$ws.Range("E30").Value = "  -0.78%  "
$ws.Range("D31").Value = "'1.850"
$ws.Range("E31").Value = "  +4.31%  "
$ws.Range("D32").Value = "'0.08897"
$ws.Range("E32").Value = "  -0.28%  "
$ws.Range("D33").Value = "'0.7641"
$ws.Range("E33").Value = "  -2.33%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'4.560"
$ws.Range("E34").Value = "  +0.99%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "'2.972"
$ws.Range("E35").Value = "  +3.34%  "
$ws.Range("D36").Value = "'1.102"
$ws.Range("E36").Value = "  +1.49%  "
$ws.Range("D37").Value = "'0.05410"
$ws.Range("E37").Value = "  -0.53%  "
$ws.Range("D38").Value = "'1.087"
$ws.Range("E38").Value = "  -0.98%  "
$ws.Range("D39").Value = "'0.01932"
$ws.Range("E39").Value = "  +0.32%  "
$ws.Range("D40").Value = "'2.817"
$ws.Range("E40").Value = "  -5.68%  "
$ws.Range("D41").Value = "'0.5076"
$ws.Range("E41").Value = "  +0.50%  "
$ws.Range("D42").Value = "'0.1653"
$ws.Range("E42").Value = "  +1.07%  "
$ws.Range("D43").Value = "'6.771"
$ws.Range("E43").Value = "  -0.38%  "
$ws.Range("D44").Value = "'8.376"
$ws.Range("E44").Value = "  +1.03%  "
$ws.Range("D45").Value = "'0.06549"
$ws.Range("E45").Value = "  -1.06%  "
$ws.Range("E46").Value = "  +0.21%  "
$ws.Range("D47").Value = "'105.16"
$ws.Range("E47").Value = "  -1.33%  "
$ws.Range("D48").Value = "'0.4671"
$ws.Range("E48").Value = "  -1.13%  "
$ws.Range("D49").Value = "'1.000"
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("D50").Value = "'1.620"
$ws.Range("E50").Value = "  -1.49%  "
$ws.Range("D51").Value = "'64.48"
$ws.Range("E51").Value = "  -0.38%  "
